$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range first so no stale cells remain
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 12338.6,            13197, 11448, 0.3305031299591065),
    @(1, 12390.96666666667,  13509, 11390, 0.3394201040267945),
    @(2, 12886.8,            13738, 12124, 0.3347977876663208),
    @(3, 12463.9,            13179, 11175, 0.324612561861674),
    @(4, 12143.5,            12875, 11001, 0.3209611415863037),
    @(5, 13056.23333333333,  13955, 11959, 0.3379259506861369),
    @(6, 12852.26666666667,  13486, 11385, 0.3383278846740723),
    @(7, 12250.03333333333,  13324, 9627,  0.3337385892868042),
    @(8, 11975.03333333333,  13058, 9878,  0.3272473971048991),
    @(9, 12213.16666666667,  13125, 11159, 0.3278856118520101)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}

$wb.Save()
